$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new structure/product ID row, matching the style of the rows above it
$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "A14-4"

$ws.Range("B6").Value = "Z45900028"

# Move selection to the next empty row, column B (matches author's post-edit state)
$ws.Range("B7").Select()
